# Add the new MSRP_VERSION column (G) to the control table.
#
# Shared-string insertion order matters (it determines the resulting
# index each new string gets in xl/sharedStrings.xml), so we write the
# data rows (G2, G3) before the header row (G1) to reproduce the target
# ordering: MSRP_2000_CHOICE, MSRP_2000_CHOICE_T, MSRP_VERSION.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G2: MSRP_2000_CHOICE (copy format from A2, a data-row cell) ---
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "MSRP_2000_CHOICE"

# --- G3: MSRP_2000_CHOICE_T (copy format from A3) ---
$ws.Range("A3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "MSRP_2000_CHOICE_T"

# --- G1: MSRP_VERSION header (copy format from A1, the header style) ---
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "MSRP_VERSION"

# New column G width (auto-fit-ish width used for the new column).
$ws.Columns.Item(7).ColumnWidth = 20

# Minor width refresh on the other columns (same visual widths, re-saved
# by a newer Excel build).
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 12.5
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 16.666666666666668

# Update the active selection to reflect where the editor left off.
$ws.Range("G8").Select()
